# daily auto push: 2026-02-16 03:18 UTC
# Insert a new data row at row 803 (pushing rows 803..844 down to 804..845)
# and populate it with the new day's entry: 2026/02/16, Mon, hour 7, rank 36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 803, shifting all rows
# from 803 downward (including the old last row 844) down by one.
$ws.Rows.Item(803).Insert()

# Column A holds date-like text (e.g. "2026/12/29"). Force text formatting
# first so Excel does not reinterpret the string as a date serial number.
$ws.Range("A803").NumberFormat = "@"
$ws.Range("A803").Value = "2026/02/16"
$ws.Range("A803").Style = "Normal"

$ws.Range("B803").Value = "月"
$ws.Range("C803").Value = 7
$ws.Range("D803").Value = 36
